$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:L1) ---
# New columns: status_seq, cr_by, cr_dtimes, upd_by, upd_dtimes, is_deleted, del_dtimes
# were inserted and the existing headers reordered.
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "descr"
$ws.Range("C1").Value = "status_seq"
$ws.Range("D1").Value = "sttyp_code"
$ws.Range("E1").Value = "lang_code"
$ws.Range("F1").Value = "is_active"
$ws.Range("G1").Value = "cr_by"
$ws.Range("H1").Value = "cr_dtimes"
$ws.Range("I1").Value = "upd_by"
$ws.Range("J1").Value = "upd_dtimes"
$ws.Range("K1").Value = "is_deleted"
$ws.Range("L1").Value = "del_dtimes"

# --- Row 2 (ACT / activated) ---
$ws.Range("A2").Value = "ACT"
$ws.Range("B2").Value = "activated"
$ws.Range("C2").Value = "NULL"
$ws.Range("D2").Value = "ADT"
$ws.Range("E2").Value = "fra"
$ws.Range("F2").Value = $true
$ws.Range("G2").Value = "superadmin"
$ws.Range("H2").Value = 45079.578406053239
$ws.Range("H2").NumberFormatLocal = "mm:ss.0"
$ws.Range("I2").Value = "NULL"
$ws.Range("J2").Value = "NULL"
$ws.Range("K2").Value = $false
$ws.Range("L2").Value = "NULL"

# --- Row 3 (DCT / disabled) ---
$ws.Range("A3").Value = "DCT"
$ws.Range("B3").Value = "disabled"
$ws.Range("C3").Value = "NULL"
$ws.Range("D3").Value = "ADT"
$ws.Range("E3").Value = "fra"
$ws.Range("F3").Value = $true
$ws.Range("G3").Value = "superadmin"
$ws.Range("H3").Value = 45079.578406053239
$ws.Range("H3").NumberFormatLocal = "mm:ss.0"
$ws.Range("I3").Value = "NULL"
$ws.Range("J3").Value = "NULL"
$ws.Range("K3").Value = $false
$ws.Range("L3").Value = "NULL"

# Move the active selection to match the author's saved cursor position.
$ws.Range("C10").Select()
